# Add 2022-Q3 data
# 1. Insert a new worksheet "2022-Q3" right after the "总计" (summary) sheet.
# 2. Populate it with the quarterly fund-holdings table (mirrors the layout
#    used by the existing quarter sheets: header row + per-fund rows).
# 3. Insert a new row at the top of the "总计" sheet's data (row 2) summarizing
#    2022-Q3, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# Use the existing "2022-Q2" sheet (header row + first data row) as a
# formatting template so fonts / borders / styles match the other quarter
# sheets exactly.
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Range("A1:H3").Copy($q3Sheet.Range("A1:H3"))

# Stamp the same row template down through row 12 (11 data rows total).
for ($r = 4; $r -le 12; $r++) {
    $q3Sheet.Range("A2:H2").Copy($q3Sheet.Range("A" + $r + ":H" + $r))
}

# ---------------------------------------------------------------------------
# Step 2: fill in the 2022-Q3 fund holdings data
# ---------------------------------------------------------------------------
$q3Data = @(
    @(0,  "010902", "博时成长领航灵活配置混合A",             "46.83", "81.70",  "5.01", "2.3462", 5),
    @(1,  "513060", "博时恒生医疗保健ETF（QDII）",           "44.00", "99.48",  "3.74", "1.6456", 6),
    @(2,  "010903", "博时成长领航灵活配置混合C",             "6.84",  "81.70",  "5.01", "0.3427", 5),
    @(3,  "513700", "鹏华中证港股通医药卫生综合ETF",         "2.96",  "94.14",  "3.05", "0.0903", 7),
    @(4,  "159892", "华夏恒生香港上市生物科技ETF（QDII）",   "1.58",  "99.13",  "4.11", "0.0649", 6),
    @(5,  "513280", "汇添富恒生香港上市生物科技ETF（QDII）", "1.51",  "100.14", "4.16", "0.0628", 6),
    @(6,  "513200", "易方达中证港股通医药卫生综合ETF",       "0.77",  "95.67",  "3.19", "0.0246", 7),
    @(7,  "159776", "银华中证港股通医药卫生综合ETF",         "0.52",  "92.74",  "3.01", "0.0157", 7),
    @(8,  "159718", "平安中证港股通医药卫生综合ETF",         "0.53",  "90.14",  "2.94", "0.0156", 7),
    @(9,  "013182", "安信港股通精选混合C",                   "0.12",  "38.51",  "1.41", "0.0017", 9),
    @(10, "013181", "安信港股通精选混合A",                   "0.02",  "38.51",  "1.41", "0.0003", 9)
)

# Columns B, D, E, F, G are stored as plain text in the source data (fund
# codes keep their leading zeros, and the percentage/NAV-like numbers keep
# their trailing zeros), so force a text number-format before assigning the
# value — otherwise Excel's COM layer auto-coerces them to floats/doubles.
$row = 2
foreach ($fund in $q3Data) {
    $q3Sheet.Cells.Item($row, 1).Value = $fund[0]

    $q3Sheet.Cells.Item($row, 2).NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 2).Value = $fund[1]

    $q3Sheet.Cells.Item($row, 3).Value = $fund[2]

    $q3Sheet.Cells.Item($row, 4).NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 4).Value = $fund[3]

    $q3Sheet.Cells.Item($row, 5).NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 5).Value = $fund[4]

    $q3Sheet.Cells.Item($row, 6).NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 6).Value = $fund[5]

    $q3Sheet.Cells.Item($row, 7).NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 7).Value = $fund[6]

    $q3Sheet.Cells.Item($row, 8).Value = $fund[7]
    $row = $row + 1
}

$q3Sheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# Step 3: insert the 2022-Q3 summary row into the "总计" sheet
# ---------------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Clear any inherited number formatting on the new row's B:D cells so they
# match the plain (unstyled) data cells used elsewhere in this column.
$summarySheet.Range("B2:D2").ClearFormats()

# Give A2 the same style used by the other index cells in column A (style
# shared with the header row).
$summarySheet.Cells.Item(1, 2).Copy()
$summarySheet.Cells.Item(2, 1).PasteSpecial(-4122)

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 11
$summarySheet.Cells.Item(2, 4).Value = 4.61

# Renumber the index column for the rows that shifted down.
for ($r = 3; $r -le 8; $r++) {
    $summarySheet.Cells.Item($r, 1).Value = $r - 2
}

$summarySheet.Range("A1").Select()
